$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous contents of the used range so we can rebuild the table
$ws.Range("A1:C5").ClearContents()

# Header row (row 1): bold/centered style 1 is already applied on A1;
# extend the same style to the new header cells B1:E1
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "2019"
$ws.Range("C1").Value = "Unnamed: 1"
$ws.Range("D1").Value = "2018"
$ws.Range("E1").Value = "Unnamed: 2"

$ws.Range("B1:E1").Style = $ws.Range("A1").Style

# Data rows
$ws.Range("B2").Value = 6.1
$ws.Range("D2").Value = 5.9

$ws.Range("B3").Value = 15.4
$ws.Range("D3").Value = 15

$ws.Range("B4").Value = 20.4
$ws.Range("D4").Value = 20.4
